$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.785.70"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.074.59"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.04"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.92"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.073.76"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.37"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000240"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.83"
$ws.Range("E14").Value = "  -3.82%  "
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "3.582.52"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "66.735.89"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.00"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.98"
$ws.Range("E19").Value = "  +3.97%  "
$ws.Range("D20").Value = "3.073.01"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "492.18"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.689"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.86"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.69"
$ws.Range("E25").Value = "  -5.79%  "
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.80"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  -5.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.64"
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").Value = "0.0₃0911"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.949"
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.99"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("E40").Value = "  -5.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.301"
$ws.Range("E41").Value = "  -3.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.32"
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("D43").Value = "2.756.51"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0346"
$ws.Range("E44").Value = "  -3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.41"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.50"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "367.14"
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.66"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.16"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("E51").Value = "  -1.85%  "
